# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51)
# to match the scraped values from the latest GitHub Actions run,
# including the Toncoin/Hedera and ImmutableX/LidoDAOToken rank swaps
# (rows 31-34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$value) {
    # Price column values such as "1.000" or "0.08230" are numeric-looking
    # text that must stay text (leading/trailing zeros matter). A plain
    # Range.Value assignment lets Excel auto-coerce those to numbers, so
    # numeric-looking strings are entered with a leading apostrophe, exactly
    # like a user typing '0.9997 into the cell.
    if ($value -match '^[+-]?\d+(\.\d+)?([eE][+-]?\d+)?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}

# Row 2
Set-TextCell $ws.Range("D2") '29.070.86'
$ws.Range("E2").Value = '  -1.71%  '

# Row 3
Set-TextCell $ws.Range("D3") '1.834.18'
$ws.Range("E3").Value = '  -1.30%  '

# Row 4
Set-TextCell $ws.Range("D4") '0.9997'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
Set-TextCell $ws.Range("D5") '240.06'
$ws.Range("E5").Value = '  -1.87%  '

# Row 6
Set-TextCell $ws.Range("D6") '0.6707'
$ws.Range("E6").Value = '  -3.34%  '

# Row 7
Set-TextCell $ws.Range("D7") '1.000'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
Set-TextCell $ws.Range("D8") '0.2969'
$ws.Range("E8").Value = '  -2.87%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.07427'
$ws.Range("E9").Value = '  -3.40%  '

# Row 10
Set-TextCell $ws.Range("D10") '22.96'
$ws.Range("E10").Value = '  -3.01%  '

# Row 11
Set-TextCell $ws.Range("D11") '0.07655'
$ws.Range("E11").Value = '  -1.38%  '

# Row 12
Set-TextCell $ws.Range("D12") '1.829.15'
$ws.Range("E12").Value = '  -1.65%  '

# Row 13
Set-TextCell $ws.Range("D13") '5.009'
$ws.Range("E13").Value = '  -2.75%  '

# Row 14
Set-TextCell $ws.Range("D14") '0.6738'
$ws.Range("E14").Value = '  -2.51%  '

# Row 15
Set-TextCell $ws.Range("D15") '86.16'
$ws.Range("E15").Value = '  -5.91%  '

# Row 16
Set-TextCell $ws.Range("D16") '6.149'
$ws.Range("E16").Value = '  -6.32%  '

# Row 17
Set-TextCell $ws.Range("D17") '29.078.08'
$ws.Range("E17").Value = '  -1.75%  '

# Row 18
Set-TextCell $ws.Range("D18") '0.000008237'
$ws.Range("E18").Value = '  -0.58%  '

# Row 19
Set-TextCell $ws.Range("D19") '227.29'
$ws.Range("E19").Value = '  -5.25%  '

# Row 20
Set-TextCell $ws.Range("D20") '12.45'
$ws.Range("E20").Value = '  -2.38%  '

# Row 21
Set-TextCell $ws.Range("D21") '0.9993'
$ws.Range("E21").Value = '  -0.13%  '

# Row 22
Set-TextCell $ws.Range("D22") '7.317'
$ws.Range("E22").Value = '  -3.68%  '

# Row 23
$ws.Range("E23").Value = '  -0.07%  '

# Row 24
Set-TextCell $ws.Range("D24") '160.37'
$ws.Range("E24").Value = '  +0.34%  '

# Row 25
Set-TextCell $ws.Range("D25") '0.1427'
$ws.Range("E25").Value = '  -4.57%  '

# Row 26
Set-TextCell $ws.Range("D26") '8.674'
$ws.Range("E26").Value = '  -2.67%  '

# Row 27
Set-TextCell $ws.Range("D27") '17.96'
$ws.Range("E27").Value = '  -1.60%  '

# Row 28
Set-TextCell $ws.Range("D28") '1.505'
$ws.Range("E28").Value = '  -1.99%  '

# Row 29
Set-TextCell $ws.Range("D29") '4.233'
$ws.Range("E29").Value = '  -0.37%  '

# Row 30
Set-TextCell $ws.Range("D30") '4.122'
$ws.Range("E30").Value = '  -1.38%  '

# Row 31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws.Range("D31") '1.194'
$ws.Range("E31").Value = '  -0.74%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws.Range("D32") '0.05395'
$ws.Range("E32").Value = '  +5.90%  '

# Row 33
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws.Range("D33") '0.7490'
$ws.Range("E33").Value = '  -2.89%  '

# Row 34
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws.Range("D34") '1.853'
$ws.Range("E34").Value = '  -1.90%  '

# Row 35
Set-TextCell $ws.Range("D35") '1.125'
$ws.Range("E35").Value = '  -2.22%  '

# Row 36
Set-TextCell $ws.Range("D36") '2.680'
$ws.Range("E36").Value = '  -0.28%  '

# Row 37
Set-TextCell $ws.Range("D37") '1.305.53'
$ws.Range("E37").Value = '  -1.77%  '

# Row 38
Set-TextCell $ws.Range("D38") '0.01805'
$ws.Range("E38").Value = '  -3.48%  '

# Row 39
Set-TextCell $ws.Range("D39") '2.710'
$ws.Range("E39").Value = '  -0.58%  '

# Row 40
Set-TextCell $ws.Range("D40") '0.9317'
$ws.Range("E40").Value = '  -3.32%  '

# Row 41
Set-TextCell $ws.Range("D41") '6.093'
$ws.Range("E41").Value = '  +5.32%  '

# Row 42
Set-TextCell $ws.Range("D42") '104.36'
$ws.Range("E42").Value = '  -1.96%  '

# Row 43
$ws.Range("E43").Value = '  -0.26%  '

# Row 44
Set-TextCell $ws.Range("D44") '0.08230'
$ws.Range("E44").Value = '  +28.70%  '

# Row 45
Set-TextCell $ws.Range("D45") '1.976.18'
$ws.Range("E45").Value = '  -1.35%  '

# Row 46
Set-TextCell $ws.Range("D46") '0.5177'
$ws.Range("E46").Value = '  -0.77%  '

# Row 47
$ws.Range("E47").Value = '  -3.33%  '

# Row 48
Set-TextCell $ws.Range("D48") '9.398'
$ws.Range("E48").Value = '  -3.44%  '

# Row 49
Set-TextCell $ws.Range("D49") '1.755'
$ws.Range("E49").Value = '  -0.88%  '

# Row 50
Set-TextCell $ws.Range("D50") '63.40'
$ws.Range("E50").Value = '  -0.14%  '

# Row 51
$ws.Range("E51").Value = '  +0.11%  '
